$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.627.19"
$ws.Range("E2").Value = "  +8.63%  "

$ws.Range("D3").Value = "3.642.01"
$ws.Range("E3").Value = "  +8.88%  "

$ws.Range("E4").Value = "  -0.08%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.69"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +7.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "193.79"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +11.76%  "

$ws.Range("E7").Value = "  +3.98%  "

$ws.Range("D8").Value = "3.602.57"
$ws.Range("E8").Value = "  +8.01%  "

$ws.Range("E9").Value = "  -0.12%  "

$ws.Range("E10").Value = "  +7.16%  "

$ws.Range("E11").Value = "  +6.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "58.30"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +10.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000299"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +8.91%  "

$ws.Range("E14").Value = "  +7.74%  "

$ws.Range("D15").Value = "4.212.64"
$ws.Range("E15").Value = "  +8.32%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "19.46"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.79%  "

$ws.Range("D17").Value = "3.628.37"
$ws.Range("E17").Value = "  +8.66%  "

$ws.Range("D18").Value = "70.503.21"
$ws.Range("E18").Value = "  +8.39%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.71"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +8.80%  "

$ws.Range("E20").Value = "  +1.22%  "

$ws.Range("E21").Value = "  +7.28%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "497.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.85%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +13.56%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.91"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +19.84%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.48"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +10.68%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "91.58"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +2.35%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +9.40%  "

$ws.Range("E28").Value = "  +8.43%  "

$ws.Range("E29").Value = "  +9.54%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "32.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.16%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.60"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +17.74%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.27"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.33%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "619.40"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.38%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "65.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.07%  "

$ws.Range("E35").Value = "  +9.80%  "

$ws.Range("D36").Value = "0.0₃0836"
$ws.Range("E36").Value = "  +14.70%  "

$ws.Range("E37").Value = "  +5.60%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.74"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.54%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "38.34"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.76%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.405"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.47%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.09%  "

$ws.Range("D42").Value = "3.337.81"
$ws.Range("E42").Value = "  +8.07%  "

$ws.Range("E43").Value = "  +12.10%  "

$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0448"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +8.72%  "

$ws.Range("B45").Value = "Fetch.AI"
$ws.Range("C45").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +12.49%  "

$ws.Range("E46").Value = "  +19.00%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.02%  "

$ws.Range("E48").Value = "  +3.95%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +8.98%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.28"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.89%  "

$ws.Range("E51").Value = "  -0.08%  "
